$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 10 (shifts existing rows 10-16 down to 12-18)
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# New row 10: Femacal de La Calera, Coquimbo, 2021-09-10 (44449), Primera
$ws.Cells.Item(10, 1).Value = 3
$ws.Cells.Item(10, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44449
$ws.Cells.Item(10, 5).Value = 5
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100107
$ws.Cells.Item(10, 8).Value = "Otros"
$ws.Cells.Item(10, 9).Value = 100107002
$ws.Cells.Item(10, 10).Value = "Chirimoya"
$ws.Cells.Item(10, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 40
$ws.Cells.Item(10, 14).Value = 30000
$ws.Cells.Item(10, 15).Value = 30000
$ws.Cells.Item(10, 16).Value = 30000
$ws.Cells.Item(10, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(10, 19).Value = 3000
$ws.Cells.Item(10, 20).Value = 10

# New row 11: Femacal de La Calera, Coquimbo, 2021-09-10 (44449), Segunda
$ws.Cells.Item(11, 1).Value = 3
$ws.Cells.Item(11, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44449
$ws.Cells.Item(11, 5).Value = 5
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100107
$ws.Cells.Item(11, 8).Value = "Otros"
$ws.Cells.Item(11, 9).Value = 100107002
$ws.Cells.Item(11, 10).Value = "Chirimoya"
$ws.Cells.Item(11, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(11, 12).Value = "Segunda"
$ws.Cells.Item(11, 13).Value = 45
$ws.Cells.Item(11, 14).Value = 27000
$ws.Cells.Item(11, 15).Value = 27000
$ws.Cells.Item(11, 16).Value = 27000
$ws.Cells.Item(11, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(11, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 19).Value = 2700
$ws.Cells.Item(11, 20).Value = 10

# Apply the same date-cell number format used by other rows in column D
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat

$wb.Save()
